$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after "OrderCreated-Event"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "PostMessage-Event"

# ---- Copy cell formatting (styles) from sheet1 so the new sheet matches ----
$ws1.Range("A1:O1").Copy() | Out-Null
$ws2.Range("A1:O1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2:O2").Copy() | Out-Null
$ws2.Range("A2:O2").PasteSpecial(-4122) | Out-Null

$ws1.Range("A3:O3").Copy() | Out-Null
$ws2.Range("A3:O3").PasteSpecial(-4122) | Out-Null

$ws1.Range("P2").Copy() | Out-Null
$ws2.Range("O2").PasteSpecial(-4122) | Out-Null

$ws1.Range("P3").Copy() | Out-Null
$ws2.Range("O3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---- Header row (row 1) ----
$ws2.Range("A1").Value2 = "TestCaseName"
$ws2.Range("B1").Value2 = "Type"
$ws2.Range("C1").Value2 = "Resource"
$ws2.Range("D1").Value2 = "TestCaseNameDesc"
$ws2.Range("E1").Value2 = "StepInfo"
$ws2.Range("F1").Value2 = "URL"
$ws2.Range("G1").Value2 = "ContentType"
$ws2.Range("H1").Value2 = "Event"
$ws2.Range("I1").Value2 = "Identifier"
$ws2.Range("J1").Value2 = "RequestContent"
$ws2.Range("K1").Value2 = "MessageType"
$ws2.Range("L1").Value2 = "Csvson"
$ws2.Range("M1").Value2 = "Action"
$ws2.Range("N1").Value2 = "StatusCode"
$ws2.Range("O1").Value2 = "Tags"

# ---- Row 2 : SendOrder ----
# NB: values are assigned in the same order the original authoring tool
# first introduced each unique string, so the generated sharedStrings.xml
# table lines up index-for-index with the target workbook.
$ws2.Range("A2").Value2 = "SendOrder"
$ws2.Range("B2").Value2 = "KAFKA"
$ws2.Range("C2").Value2 = "json"
$ws2.Range("D2").Value2 = "Create order "
$ws2.Range("E2").Value2 = "As a user needs to;create order;details;order;user;"
$ws2.Range("H2").Value2 = "OrderCreated"
$ws2.Range("O2").Value2 = " @simple-send"

$json2 = @"
{
  "customer": {
    "customerId": 1001,
    "firstname": "Ronnie",
    "lastname": "Sander"
  },
  "orderNumber" : "1234-1234-1234",
  "orderDesc": "Order Dell Laptop",
  "orderStatus": "Started",
  "purchasedProducts": [
    {
      "productId": 901,
      "productName": "Dell Inspiron 3583 15",
      "productDesc": "Laptop Intel Celeron – 128GB SSD – 4GB DDR4 – 1.6GHz - Intel UHD Graphics 610 - Windows 10 Home in S Mode - Inspiron 15 3000 Series"
    }
  ]
}
"@
$ws2.Range("J2").Value2 = $json2

$ws2.Range("K2").Value2 = "JSONType"

# ---- Row 3 : VERIFY_ORDER_CREATED_EVENT ----
$ws2.Range("A3").Value2 = "VERIFY_ORDER_CREATED_EVENT"
$ws2.Range("B3").Value2 = "KAFKA"
$ws2.Range("C3").Value2 = "json"
$ws2.Range("D3").Value2 = "Validate created order event"
$ws2.Range("E3").Value2 = "contains order information"
$ws2.Range("H3").Value2 = "OrderCreated"
$ws2.Range("I3").Value2 = "1234-1234-1234"
$ws2.Range("K3").Value2 = "JSONType"

$csv3 = @"
orderNumber,orderDesc,orderStatus, customer/customerId:firstname:lastname,orderStatus
1234-1234-1234,Order Dell Laptop,Started,i~1001:Ronnie:Sander,Started
"@
$ws2.Range("L3").Value2 = $csv3

$ws2.Range("O3").Value2 = "  @validate_kafka_message @IncludesByPath"

# ---- Row heights / column widths (match target layout as closely as possible) ----
$ws2.Rows.Item(1).RowHeight = 15.75
$ws2.Rows.Item(2).RowHeight = 299.25
$ws2.Rows.Item(3).RowHeight = 31.5

$ws2.Columns.Item(1).ColumnWidth = 33.140625
$ws2.Columns.Item(3).ColumnWidth = 12.42578125
$ws2.Columns.Item(4).ColumnWidth = 28
$ws2.Columns.Item(5).ColumnWidth = 48.7109375
$ws2.Columns.Item(6).ColumnWidth = 24.5703125
$ws2.Columns.Item(7).ColumnWidth = 16.140625
$ws2.Columns.Item(8).ColumnWidth = 13.7109375
$ws2.Columns.Item(9).ColumnWidth = 16.28515625
$ws2.Columns.Item(10).ColumnWidth = 72.5703125
$ws2.Columns.Item(11).ColumnWidth = 14
$ws2.Columns.Item(12).ColumnWidth = 91
$ws2.Columns.Item(15).ColumnWidth = 45

# ---- Tidy up the selection left on the original sheet ----
$ws1.Activate()
$ws1.Range("I1").Select() | Out-Null

# ---- Make the new sheet the active / selected sheet ----
$ws2.Activate()
$ws2.Range("E2").Select() | Out-Null

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host " -" $s.Name
}
